$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header: "pvalues" -> "p"
$ws.Range("F1").Value = "p"

# F2, F3, F4 p-values corrected to ".401" (kept as text, not numeric,
# to match the original inline-string cell type / formatting).
$ws.Range("F2:F4").Formula = '=".401"'
$ws.Range("F2:F4").Copy()
$ws.Range("F2:F4").PasteSpecial(-4163)  # xlPasteValues

# F6 p-value corrected to ".519" (also kept as text).
$ws.Range("F6").Formula = '=".519"'
$ws.Range("F6").Copy()
$ws.Range("F6").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0

# New footnote row describing the multiple-testing correction applied.
$ws.Range("A9").Value = "Multiple tests correction applied to p values: Benjamini-Hochberg"
